$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.650.70"
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("D3").Value = "3.500.01"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.98"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.20"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E8").Value = "  +4.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.132"
$ws.Range("E9").Value = "  +7.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.34"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.433"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "4.105.39"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.30"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000181"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("D16").Value = "66.691.60"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "3.498.33"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.04"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.92"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.96"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.34"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.535"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.07"
$ws.Range("E31").Value = "  +1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.80"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.40"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.61"
$ws.Range("E34").Value = "  +5.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "162.19"
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.903"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.80"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.68"
$ws.Range("E39").Value = "  +4.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0745"
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.52"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.27"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("D43").Value = "2.809.81"
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.95"
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.56"
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0313"
$ws.Range("E46").Value = "  +1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "343.01"
$ws.Range("E47").Value = "  +0.30%  "
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.13"
$ws.Range("E49").Value = "  +5.18%  "
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.51"
$ws.Range("E51").Value = "  +1.88%  "
